$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the header row: "<col>_old" -> "<col>_FV2404", "<col>_new" -> "<col>_FV2410"
# ---------------------------------------------------------------------------
$headerMap = @{
    "A1" = "Segmentname_FV2404"
    "B1" = "Segmentgruppe_FV2404"
    "C1" = "Segment_FV2404"
    "D1" = "Datenelement_FV2404"
    "E1" = "Segment ID_FV2404"
    "F1" = "Code_FV2404"
    "G1" = "Qualifier_FV2404"
    "H1" = "Beschreibung_FV2404"
    "I1" = "Bedingungsausdruck_FV2404"
    "J1" = "Bedingung_FV2404"
    "L1" = "Segmentname_FV2410"
    "M1" = "Segmentgruppe_FV2410"
    "N1" = "Segment_FV2410"
    "O1" = "Datenelement_FV2410"
    "P1" = "Segment ID_FV2410"
    "Q1" = "Code_FV2410"
    "R1" = "Qualifier_FV2410"
    "S1" = "Beschreibung_FV2410"
    "T1" = "Bedingungsausdruck_FV2410"
    "U1" = "Bedingung_FV2410"
}

foreach ($addr in $headerMap.Keys) {
    $ws.Range($addr).Value = $headerMap[$addr]
}

# ---------------------------------------------------------------------------
# 2. Turn the header + data range into a real Excel Table ("Table1"), with
#    column headers picked up from row 1 (now already renamed above).
#
#    The header row already carries bold/fill/border direct formatting
#    (no dxf yet). Adding a ListObject over already-formatted header cells
#    makes Excel snapshot that look into a new dxf (headerRowDxfId) AND
#    assign a named TableStyle ("TableStyleMedium2") to tableStyleInfo -
#    neither of which is present in the target workbook (plain/no table
#    style, no header dxf). To reproduce that, we park a copy of the
#    header formatting on a scratch range, blank the header formatting
#    back to the default "Normal" style, create the table, then paste the
#    formatting back in a single shot (so the final style maps onto the
#    very same cellXf that was already in the file, instead of leaving
#    behind extra half-applied style entries).
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$donor = $ws.Range("A1000:U1000")

$headerRange.Copy()
$donor.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$headerRange.Style = "Normal"

$tableRange = $ws.Range("A1:U67")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$donor.Copy()
$headerRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$donor.Clear()

# ---------------------------------------------------------------------------
# 3. Freeze the header row (split/freeze pane below row 1).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
